$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.398.09"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.821.88"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.41"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5139"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3845"
$ws.Range("E8").Value = "  -1.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08229"
$ws.Range("E9").Value = "  +6.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.118"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.82"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.346"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.03"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.004"
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.435"
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").Value = "1.819.87"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.90"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001109"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06624"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.037"
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("D23").Value = "28.431.50"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.49"
$ws.Range("E24").Value = "  +2.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.247"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.86"
$ws.Range("E26").Value = "  +1.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.95"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("D28").Value = "2.028.09"
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.395"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.08"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.702"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.683"
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07423"
$ws.Range("E35").Value = "  +5.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.36"
$ws.Range("E36").Value = "  +9.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2214"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02341"
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.183"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.749"
$ws.Range("E40").Value = "  -2.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6341"
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.387"
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.60"
$ws.Range("E44").Value = "  +1.46%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6156"
$ws.Range("E45").Value = "  +4.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.814"
$ws.Range("E46").Value = "  +2.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.25"
$ws.Range("E47").Value = "  +1.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.998"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.202"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06912"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.072"
$ws.Range("E51").Value = "  +1.02%  "
